# Comparison.xlsx - "hopeful simple fix overlook"
# Adds the actual-steps timing data (C/E numeric columns + G "Xms" labels)
# that was missing for the 20x20 / 20x50 maze rows, fixes the
# "DFS actual stesp" -> "DFS actual steps" typo, and moves the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) --------------------------------------------------
# D1/E1 swap meaning: D1 becomes "BFS actual steps", E1 becomes the
# (typo-corrected) "DFS actual steps". F1/G1 stay "BFS time"/"DFS time".
$ws.Range("D1").Value = "BFS actual steps"
$ws.Range("E1").Value = "DFS actual steps"
$ws.Range("F1").Value = "BFS time"
$ws.Range("G1").Value = "DFS time"
$ws.Range("H1").Value = "File Name"

# --- Newly-filled data rows (2-7) ----------------------------------------
$ws.Range("C2").Value = 92
$ws.Range("E2").Value = 294
$ws.Range("G2").Value = "34ms"

$ws.Range("C3").Value = 4
$ws.Range("E3").Value = 4
$ws.Range("G3").Value = "32ms"

$ws.Range("C4").Value = 240
$ws.Range("E4").Value = 331
$ws.Range("G4").Value = "34ms"

$ws.Range("C5").Value = 185
$ws.Range("E5").Value = 348
$ws.Range("G5").Value = "34ms"

$ws.Range("C6").Value = 118
$ws.Range("E6").Value = 390
$ws.Range("G6").Value = "35ms"

$ws.Range("C7").Value = 447
$ws.Range("E7").Value = 780
$ws.Range("G7").Value = "35ms"

# --- Selection moved from H17 to F14 -------------------------------------
$ws.Range("F14").Select() | Out-Null
